{"js": "// Update the answers of the two-digit/one-digit division problems.\n// Each cell in the table holds one \"A\u00f7B=C, D\" style answer; find the\n// old answer text in the body and replace it with the new one.\n\nconst replacements = [\n  [\"19\u00f77=2, 5\", \"27\u00f73=9, 0\"],\n  [\"54\u00f72=27, 0\", \"93\u00f75=18, 3\"],\n  [\"60\u00f79=6, 6\", \"72\u00f78=9, 0\"],\n  [\"89\u00f77=12, 5\", \"17\u00f76=2, 5\"],\n  [\"66\u00f72=33, 0\", \"49\u00f78=6, 1\"],\n  [\"58\u00f74=14, 2\", \"77\u00f76=12, 5\"],\n  [\"18\u00f77=2, 4\", \"71\u00f75=14, 1\"],\n  [\"74\u00f72=37, 0\", \"69\u00f73=23, 0\"],\n  [\"86\u00f76=14, 2\", \"39\u00f78=4, 7\"],\n  [\"72\u00f75=14, 2\", \"18\u00f78=2, 2\"],\n  [\"80\u00f78=10, 0\", \"33\u00f74=8, 1\"],\n  [\"98\u00f78=12, 2\", \"69\u00f76=11, 3\"],\n  [\"59\u00f79=6, 5\", \"58\u00f72=29, 0\"],\n  [\"13\u00f74=3, 1\", \"79\u00f78=9, 7\"],\n  [\"94\u00f73=31, 1\", \"77\u00f72=38, 1\"],\n  [\"21\u00f75=4, 1\", \"65\u00f74=16, 1\"],\n  [\"51\u00f77=7, 2\", \"30\u00f74=7, 2\"],\n  [\"27\u00f77=3, 6\", \"65\u00f77=9, 2\"],\n  [\"99\u00f77=14, 1\", \"81\u00f78=10, 1\"],\n  [\"42\u00f76=7, 0\", \"65\u00f73=21, 2\"],\n  [\"87\u00f73=29, 0\", \"99\u00f75=19, 4\"],\n  [\"61\u00f76=10, 1\", \"83\u00f74=20, 3\"],\n  [\"49\u00f73=16, 1\", \"17\u00f79=1, 8\"],\n  [\"53\u00f78=6, 5\", \"26\u00f75=5, 1\"],\n  [\"29\u00f75=5, 4\", \"64\u00f79=7, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the answers of the two-digit/one-digit division problems.\n# Each cell in the table holds one \"A\u00f7B=C, D\" style answer; replace the\n# old answer text with the new one, cell by cell, using Find/Replace so\n# we only ever touch the exact run that currently holds that text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"19\u00f77=2, 5\", \"27\u00f73=9, 0\"),\n    @(\"54\u00f72=27, 0\", \"93\u00f75=18, 3\"),\n    @(\"60\u00f79=6, 6\", \"72\u00f78=9, 0\"),\n    @(\"89\u00f77=12, 5\", \"17\u00f76=2, 5\"),\n    @(\"66\u00f72=33, 0\", \"49\u00f78=6, 1\"),\n    @(\"58\u00f74=14, 2\", \"77\u00f76=12, 5\"),\n    @(\"18\u00f77=2, 4\", \"71\u00f75=14, 1\"),\n    @(\"74\u00f72=37, 0\", \"69\u00f73=23, 0\"),\n    @(\"86\u00f76=14, 2\", \"39\u00f78=4, 7\"),\n    @(\"72\u00f75=14, 2\", \"18\u00f78=2, 2\"),\n    @(\"80\u00f78=10, 0\", \"33\u00f74=8, 1\"),\n    @(\"98\u00f78=12, 2\", \"69\u00f76=11, 3\"),\n    @(\"59\u00f79=6, 5\", \"58\u00f72=29, 0\"),\n    @(\"13\u00f74=3, 1\", \"79\u00f78=9, 7\"),\n    @(\"94\u00f73=31, 1\", \"77\u00f72=38, 1\"),\n    @(\"21\u00f75=4, 1\", \"65\u00f74=16, 1\"),\n    @(\"51\u00f77=7, 2\", \"30\u00f74=7, 2\"),\n    @(\"27\u00f77=3, 6\", \"65\u00f77=9, 2\"),\n    @(\"99\u00f77=14, 1\", \"81\u00f78=10, 1\"),\n    @(\"42\u00f76=7, 0\", \"65\u00f73=21, 2\"),\n    @(\"87\u00f73=29, 0\", \"99\u00f75=19, 4\"),\n    @(\"61\u00f76=10, 1\", \"83\u00f74=20, 3\"),\n    @(\"49\u00f73=16, 1\", \"17\u00f79=1, 8\"),\n    @(\"53\u00f78=6, 5\", \"26\u00f75=5, 1\"),\n    @(\"29\u00f75=5, 4\", \"64\u00f79=7, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\n$d.Save()\n"}
